$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$style = $ws.Range('D2').Style
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '63.606.68'
$ws.Range('D2').Style = $style
$ws.Range('E2').Value = '  +4.57%  '

# Row 3
$style = $ws.Range('D3').Style
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.427.65'
$ws.Range('D3').Style = $style
$ws.Range('E3').Value = '  +5.87%  '

# Row 4
$ws.Range('E4').Value = '  +0.02%  '

# Row 5
$style = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '576.61'
$ws.Range('D5').Style = $style
$ws.Range('E5').Value = '  +6.51%  '

# Row 6
$style = $ws.Range('D6').Style
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '156.89'
$ws.Range('D6').Style = $style
$ws.Range('E6').Value = '  +6.39%  '

# Row 7
$style = $ws.Range('D7').Style
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').Style = $style
$ws.Range('E7').Value = '  +0.02%  '

# Row 8
$style = $ws.Range('D8').Style
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.430.64'
$ws.Range('D8').Style = $style
$ws.Range('E8').Value = '  +5.75%  '

# Row 9
$ws.Range('E9').Value = '  +1.22%  '

# Row 10
$ws.Range('E10').Value = '  +2.87%  '

# Row 11
$ws.Range('E11').Value = '  +7.14%  '

# Row 12
$ws.Range('E12').Value = '  +0.31%  '

# Row 13
$style = $ws.Range('D13').Style
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.019.77'
$ws.Range('D13').Style = $style
$ws.Range('E13').Value = '  +6.01%  '

# Row 14
$ws.Range('E14').Value = '  -0.64%  '

# Row 15
$ws.Range('E15').Value = '  +6.88%  '

# Row 16
$style = $ws.Range('D16').Style
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '27.39'
$ws.Range('D16').Style = $style
$ws.Range('E16').Value = '  +4.54%  '

# Row 17
$style = $ws.Range('D17').Style
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '63.661.74'
$ws.Range('D17').Style = $style
$ws.Range('E17').Value = '  +4.72%  '

# Row 18
$style = $ws.Range('D18').Style
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.438.18'
$ws.Range('D18').Style = $style
$ws.Range('E18').Value = '  +5.96%  '

# Row 19
$style = $ws.Range('D19').Style
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.43'
$ws.Range('D19').Style = $style
$ws.Range('E19').Value = '  +1.79%  '

# Row 20
$style = $ws.Range('D20').Style
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.36'
$ws.Range('D20').Style = $style
$ws.Range('E20').Value = '  +7.62%  '

# Row 21
$style = $ws.Range('D21').Style
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '8.50'
$ws.Range('D21').Style = $style
$ws.Range('E21').Value = '  +1.67%  '

# Row 22
$style = $ws.Range('D22').Style
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '392.25'
$ws.Range('D22').Style = $style
$ws.Range('E22').Value = '  +3.99%  '

# Row 23
$style = $ws.Range('D23').Style
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.00'
$ws.Range('D23').Style = $style
$ws.Range('E23').Value = '  -0.24%  '

# Row 24
$style = $ws.Range('D24').Style
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.539'
$ws.Range('D24').Style = $style
$ws.Range('E24').Value = '  +2.16%  '

# Row 25
$style = $ws.Range('D25').Style
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '71.98'
$ws.Range('D25').Style = $style
$ws.Range('E25').Value = '  +2.79%  '

# Row 26
$ws.Range('E26').Value = '  +18.97%  '

# Row 27
$style = $ws.Range('D27').Style
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.56'
$ws.Range('D27').Style = $style
$ws.Range('E27').Value = '  +10.09%  '

# Row 28
$ws.Range('E28').Value = '  +5.04%  '

# Row 29
$ws.Range('E29').Value = '  -0.03%  '

# Row 30
$style = $ws.Range('D30').Style
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.72'
$ws.Range('D30').Style = $style
$ws.Range('E30').Value = '  +8.05%  '

# Row 31
$style = $ws.Range('D31').Style
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.38'
$ws.Range('D31').Style = $style
$ws.Range('E31').Value = '  +12.95%  '

# Row 32
$ws.Range('E32').Value = '  +6.56%  '

# Row 33
$style = $ws.Range('D33').Style
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.84'
$ws.Range('D33').Style = $style
$ws.Range('E33').Value = '  +8.91%  '

# Row 34
$style = $ws.Range('D34').Style
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '23.58'
$ws.Range('D34').Style = $style
$ws.Range('E34').Value = '  +4.47%  '

# Row 35
$style = $ws.Range('D35').Style
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.998'
$ws.Range('D35').Style = $style
$ws.Range('E35').Value = '  -0.08%  '

# Row 36
$style = $ws.Range('D36').Style
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.85'
$ws.Range('D36').Style = $style
$ws.Range('E36').Value = '  +3.53%  '

# Row 37
$ws.Range('E37').Value = '  +4.75%  '

# Row 38
$style = $ws.Range('D38').Style
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '157.64'
$ws.Range('D38').Style = $style
$ws.Range('E38').Value = '  -0.35%  '

# Row 39
$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$style = $ws.Range('D39').Style
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0786'
$ws.Range('D39').Style = $style
$ws.Range('E39').Value = '  +10.10%  '

# Row 40
$ws.Range('B40').Value = 'EnergySwap'
$ws.Range('C40').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$style = $ws.Range('D40').Style
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '28.12'
$ws.Range('D40').Style = $style
$ws.Range('E40').Value = '  +6.20%  '

# Row 41
$ws.Range('E41').Value = '  +8.30%  '

# Row 42
$style = $ws.Range('D42').Style
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.866.27'
$ws.Range('D42').Style = $style
$ws.Range('E42').Value = '  +1.95%  '

# Row 43
$ws.Range('E43').Value = '  +1.41%  '

# Row 44
$style = $ws.Range('D44').Style
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '41.96'
$ws.Range('D44').Style = $style
$ws.Range('E44').Value = '  +4.96%  '

# Row 45
$style = $ws.Range('D45').Style
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.769'
$ws.Range('D45').Style = $style
$ws.Range('E45').Value = '  +6.12%  '

# Row 46
$style = $ws.Range('D46').Style
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.40'
$ws.Range('D46').Style = $style
$ws.Range('E46').Value = '  +3.19%  '

# Row 47
$ws.Range('E47').Value = '  +9.28%  '

# Row 48
$style = $ws.Range('D48').Style
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.475.43'
$ws.Range('D48').Style = $style
$ws.Range('E48').Value = '  +6.04%  '

# Row 49
$style = $ws.Range('D49').Style
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '22.59'
$ws.Range('D49').Style = $style
$ws.Range('E49').Value = '  +7.23%  '

# Row 50
$ws.Range('E50').Value = '  +22.87%  '

# Row 51
$ws.Range('E51').Value = '  +2.72%  '
